$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) to the s_vals sheet, matching the existing
# header formatting (bold, centered, bordered) used by the other headers.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new column's data rows with the default value of 0.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
